# Update "想去人数" (F) counts and mark one ticket tier "不可售" (G3)
# on both the "展览" and "全部类型" worksheets. The two sheets list largely
# the same events, but "全部类型" has one extra row (a concert) inserted
# at row 6, so the row numbers differ from row 6 onward.

$wb = $excel.ActiveWorkbook

# Row -> column -> new value, for the "展览" worksheet.
$updatesExhibition = @{
    2  = @{ F = 172 }
    3  = @{ G = "不可售" }
    4  = @{ F = 582 }
    5  = @{ F = 1847 }
    6  = @{ F = 480 }
    9  = @{ F = 2466 }
    10 = @{ F = 150 }
    11 = @{ F = 82 }
    12 = @{ F = 169 }
    13 = @{ F = 1493 }
    21 = @{ F = 216 }
    23 = @{ F = 6 }
    24 = @{ F = 143 }
    26 = @{ F = 1558 }
    27 = @{ F = 23 }
    28 = @{ F = 387 }
    29 = @{ F = 432 }
    30 = @{ F = 197 }
    32 = @{ F = 397 }
}

# Row -> column -> new value, for the "全部类型" worksheet (offset by +1
# from row 6 onward because of the extra concert row at row 6).
$updatesAllTypes = @{
    2  = @{ F = 172 }
    3  = @{ G = "不可售" }
    4  = @{ F = 582 }
    5  = @{ F = 1847 }
    7  = @{ F = 480 }
    10 = @{ F = 2466 }
    11 = @{ F = 150 }
    12 = @{ F = 82 }
    13 = @{ F = 169 }
    14 = @{ F = 1493 }
    22 = @{ F = 216 }
    24 = @{ F = 6 }
    25 = @{ F = 143 }
    27 = @{ F = 1558 }
    28 = @{ F = 23 }
    29 = @{ F = 387 }
    30 = @{ F = 432 }
    31 = @{ F = 197 }
    33 = @{ F = 397 }
}

$sheetUpdates = @{
    "展览"   = $updatesExhibition
    "全部类型" = $updatesAllTypes
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $updates = $sheetUpdates[$sheetName]

    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$row").Value = $cols[$col]
        }
    }
}
